$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Paragraph 1: "Basic if demonstration :" ---
# Add _GoBack bookmark + spellStart/spellEnd proofErr markers at the start,
# and split the trailing " demonstration :" run into " " / "demonstration" / " :"
$p1 = $d.Paragraphs(1).Range
$xml1 = @"
<w:p xmlns:w="$wNs" w:rsidP="00F5495F" w:rsidR="00BA34AE" w:rsidRDefault="00BA34AE"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">Basic </w:t></w:r><w:r w:rsidR="002A1F2A"><w:t>if</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>demonstration</w:t></w:r><w:r><w:t> :</w:t></w:r></w:p>
"@
[void]$p1.InsertXML($xml1)

# --- Paragraph 2: "The ELSEIF paragraph." ---
# Add spellStart/spellEnd proofErr markers at the start, merge "The "+"ELSEIF" into
# "The ELSEIF " and split the trailing " paragraph." run into "paragraph" / "."
$p2 = $d.Paragraphs(2).Range
$xml2 = @"
<w:p xmlns:w="$wNs" w:rsidP="00430772" w:rsidR="00430772" w:rsidRDefault="00430772"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">The ELSEIF </w:t></w:r><w:r><w:t>paragraph</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
"@
[void]$p2.InsertXML($xml2)

# --- Paragraph 3: "End of demonstration." ---
# Add spellStart/spellEnd proofErr markers at the start, and split the trailing
# " of demonstration." run into " of " / "demonstration" / "."
$p3 = $d.Paragraphs(3).Range
$xml3 = @"
<w:p xmlns:w="$wNs" w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00BA34AE"><w:proofErr w:type="spellStart"/><w:proofErr w:type="spellEnd"/><w:r><w:t>En</w:t></w:r><w:r w:rsidR="006B5B12"><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:r><w:t>demonstration</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
"@
[void]$p3.InsertXML($xml3)
